# Mercury Tours - datosBancos.xlsx
# "Crear cuenta bancaria terminada"
#
# Insert a new row (the "codigo" row) above the existing data row, listing
# the numeric bank codes (01-14) that line up under each bank name in row 2,
# plus a trailing helper value (15) in column P.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing row 3 ("Incorrecto", 123, "banco*") down to row 4 and
# open up a fresh row 3 for the new "codigo" data.
$ws.Rows("3:3").Insert()

# Give the new column used by the codigo row (D) the same kind of explicit
# width the other data columns already have.
$ws.Columns("D:D").ColumnWidth = 11.5

# Row 3: "codigo" label plus the bank codes 01..14 across columns B..O.
$ws.Range("A3").Value = "codigo"

# Bank codes must be stored as text (e.g. "01", not 1), so format as text
# before writing the values.
$ws.Range("B3:O3").NumberFormat = "@"

$codigos = @("01","02","03","04","05","06","07","08","09","10","11","12","13","14")
for ($i = 0; $i -lt $codigos.Length; $i++) {
    $ws.Cells.Item(3, 2 + $i).Value = $codigos[$i]
}

# Trailing numeric marker in column P.
$ws.Range("P3").Value = 15

# Leave the workbook scrolled/selected near the newly added data.
$ws.Range("B11").Select() | Out-Null
